$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034898999729031
$ws.Range("D2").Value = 1.037348674143905
$ws.Range("E2").Value = 1.042813083993994
$ws.Range("F2").Value = 1.050877323449932
$ws.Range("I2").Value = 1.036562890441823
$ws.Range("J2").Value = 1.040015721869943
$ws.Range("K2").Value = 1.040139813858282
$ws.Range("L2").Value = 1.045588712784967
$ws.Range("M2").Value = 1.053630381151885
$ws.Range("N2").Value = 1.041492663486469
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03585520950967
$ws.Range("D3").Value = 1.037865259363114
$ws.Range("E3").Value = 1.043698400329391
$ws.Range("F3").Value = 1.051927246086503
$ws.Range("I3").Value = 1.036728675710002
$ws.Range("J3").Value = 1.040615221469768
$ws.Range("K3").Value = 1.040466965015093
$ws.Range("L3").Value = 1.046284742189328
$ws.Range("M3").Value = 1.054492226830413
$ws.Range("N3").Value = 1.04209301444449
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036474369398514
$ws.Range("D4").Value = 1.038199658230909
$ws.Range("E4").Value = 1.04427206426585
$ws.Range("F4").Value = 1.052607761844343
$ws.Range("I4").Value = 1.036834864325452
$ws.Range("J4").Value = 1.041002955450372
$ws.Range("K4").Value = 1.040678072916442
$ws.Range("L4").Value = 1.046735288193394
$ws.Range("M4").Value = 1.055050440215963
$ws.Range("N4").Value = 1.04248129905182
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036734765925344
$ws.Range("D5").Value = 1.038340269785279
$ws.Range("E5").Value = 1.044513424180367
$ws.Range("F5").Value = 1.052894123879218
$ws.Range("I5").Value = 1.036879245543404
$ws.Range("J5").Value = 1.041165914326828
$ws.Range("K5").Value = 1.040766682485641
$ws.Range("L5").Value = 1.046924737149593
$ws.Range("M5").Value = 1.055285241814959
$ws.Range("N5").Value = 1.042644489348572
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036778493588944
$ws.Range("D6").Value = 1.038363880808495
$ws.Range("E6").Value = 1.044553960774239
$ws.Range("F6").Value = 1.052942221315683
$ws.Range("I6").Value = 1.036886682053445
$ws.Range("J6").Value = 1.041193273206167
$ws.Range("K6").Value = 1.040781552167646
$ws.Range("L6").Value = 1.046956548747984
$ws.Range("M6").Value = 1.055324673555065
$ws.Range("N6").Value = 1.042671887080658
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036477848431001
$ws.Range("D7").Value = 1.038201536972631
$ws.Range("E7").Value = 1.044275288577427
$ws.Range("F7").Value = 1.052611587155648
$ws.Range("I7").Value = 1.036835458373695
$ws.Range("J7").Value = 1.04100513309071
$ws.Range("K7").Value = 1.040679257474158
$ws.Range("L7").Value = 1.046737819465877
$ws.Range("M7").Value = 1.05505357714196
$ws.Range("N7").Value = 1.042483479784656
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035222066148574
$ws.Range("D8").Value = 1.037523227815018
$ws.Range("E8").Value = 1.043112113414372
$ws.Range("F8").Value = 1.051231911882934
$ws.Range("I8").Value = 1.036619142758973
$ws.Range("J8").Value = 1.040218362820618
$ws.Range("K8").Value = 1.040250495683082
$ws.Range("L8").Value = 1.045823903758314
$ws.Range("M8").Value = 1.05392153320848
$ws.Range("N8").Value = 1.041695592210537
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033012532618427
$ws.Range("D9").Value = 1.036329060810094
$ws.Range("E9").Value = 1.041068666320975
$ws.Range("F9").Value = 1.048809568133686
$ws.Range("I9").Value = 1.036229679406873
$ws.Range("J9").Value = 1.038830613344115
$ws.Range("K9").Value = 1.039490565522651
$ws.Range("L9").Value = 1.044214797631286
$ws.Range("M9").Value = 1.051930917967921
$ws.Range("N9").Value = 1.040305871970601
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031541785847536
$ws.Range("D10").Value = 1.035533792507203
$ws.Range("E10").Value = 1.039710614519688
$ws.Range("F10").Value = 1.047200663891946
$ws.Range("I10").Value = 1.035964500351496
$ws.Range("J10").Value = 1.037904581644352
$ws.Range("K10").Value = 1.038981057613389
$ws.Range("L10").Value = 1.043143007682289
$ws.Range("M10").Value = 1.050606715410585
$ws.Range("N10").Value = 1.039378525199606
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030905485323439
$ws.Range("D11").Value = 1.035189654297265
$ws.Range("E11").Value = 1.03912358404848
$ws.Range("F11").Value = 1.046505424232113
$ws.Range("I11").Value = 1.035848368994352
$ws.Range("J11").Value = 1.037503404458687
$ws.Range("K11").Value = 1.038759764540901
$ws.Range("L11").Value = 1.0426791464524
$ws.Range("M11").Value = 1.050034014405235
$ws.Range("N11").Value = 1.038976778296323
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030669217175818
$ws.Range("D12").Value = 1.035061860590392
$ws.Range("E12").Value = 1.038905688352399
$ws.Range("F12").Value = 1.046247396545077
$ws.Range("I12").Value = 1.035805036833723
$ws.Range("J12").Value = 1.037354360207613
$ws.Range("K12").Value = 1.038677466514745
$ws.Range("L12").Value = 1.042506883269259
$ws.Range("M12").Value = 1.049821391922462
$ws.Range("N12").Value = 1.038827522385316
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030719893764749
$ws.Range("D13").Value = 1.035089271182755
$ws.Range("E13").Value = 1.038952420752698
$ws.Range("F13").Value = 1.046302734602534
$ws.Range("I13").Value = 1.035814340580654
$ws.Range("J13").Value = 1.0373863320136
$ws.Range("K13").Value = 1.038695124224108
$ws.Range("L13").Value = 1.042543832685399
$ws.Range("M13").Value = 1.049866995438536
$ws.Range("N13").Value = 1.038859539594935
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030885953643827
$ws.Range("D14").Value = 1.035179090114668
$ws.Range("E14").Value = 1.03910556958348
$ws.Range("F14").Value = 1.046484091185811
$ws.Range("I14").Value = 1.035844791135478
$ws.Range("J14").Value = 1.037491085005251
$ws.Range("K14").Value = 1.03875296379219
$ws.Range("L14").Value = 1.04266490639424
$ws.Range("M14").Value = 1.050016436834282
$ws.Range("N14").Value = 1.038964441347849
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030988279528775
$ws.Range("D15").Value = 1.035234435153459
$ws.Range("E15").Value = 1.039199950012949
$ws.Range("F15").Value = 1.046595859517642
$ws.Range("I15").Value = 1.035863526798446
$ws.Range("J15").Value = 1.03755562293682
$ws.Range("K15").Value = 1.038788587452297
$ws.Range("L15").Value = 1.042739508635645
$ws.Range("M15").Value = 1.05010852644063
$ws.Range("N15").Value = 1.039029070930684
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031584026506932
$ws.Range("D16").Value = 1.035556636574941
$ws.Range("E16").Value = 1.03974959527391
$ws.Range("F16").Value = 1.047246834795187
$ws.Range("I16").Value = 1.035972180117204
$ws.Range("J16").Value = 1.037931202318372
$ws.Range("K16").Value = 1.038995730029417
$ws.Range("L16").Value = 1.043173797560873
$ws.Range("M16").Value = 1.050644738249625
$ws.Range("N16").Value = 1.039405183678038
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031957868583021
$ws.Range("D17").Value = 1.035758804682621
$ws.Range("E17").Value = 1.040094645708189
$ws.Range("F17").Value = 1.047655556980085
$ws.Range("I17").Value = 1.036039985933867
$ws.Range("J17").Value = 1.038166740511924
$ws.Range("K17").Value = 1.039125485786056
$ws.Range("L17").Value = 1.043446277939838
$ws.Range("M17").Value = 1.050981274425363
$ws.Range("N17").Value = 1.039641056362841
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032175976641182
$ws.Range("D18").Value = 1.03587674685594
$ws.Range("E18").Value = 1.040296005641411
$ws.Range("F18").Value = 1.047894095286627
$ws.Range("I18").Value = 1.036079409767388
$ws.Range("J18").Value = 1.038304106669337
$ws.Range("K18").Value = 1.039201105132422
$ws.Range("L18").Value = 1.043605233355872
$ws.Range("M18").Value = 1.051177636551071
$ws.Range("N18").Value = 1.039778617595953
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032250354702725
$ws.Range("D19").Value = 1.035916965586611
$ws.Range("E19").Value = 1.040364680739779
$ws.Range("F19").Value = 1.047975453989975
$ws.Range("I19").Value = 1.036092830860907
$ws.Range("J19").Value = 1.038350941684602
$ws.Range("K19").Value = 1.039226878316619
$ws.Range("L19").Value = 1.043659436788438
$ws.Range("M19").Value = 1.051244602174691
$ws.Range("N19").Value = 1.039825519122312
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031917753451594
$ws.Range("D20").Value = 1.035737111779672
$ws.Range("E20").Value = 1.040057614932935
$ws.Range("F20").Value = 1.047611690721126
$ws.Range("I20").Value = 1.036032724055839
$ws.Range("J20").Value = 1.038141471482122
$ws.Range("K20").Value = 1.039111570935382
$ws.Range("L20").Value = 1.043417041069369
$ws.Range("M20").Value = 1.050945160417764
$ws.Range("N20").Value = 1.039615751448118
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030837050873736
$ws.Range("D21").Value = 1.035152639710103
$ws.Range("E21").Value = 1.039060466822314
$ws.Range("F21").Value = 1.046430680241737
$ws.Range("I21").Value = 1.035835829606608
$ws.Range("J21").Value = 1.037460238659493
$ws.Range("K21").Value = 1.038735934230207
$ws.Range("L21").Value = 1.04262925221706
$ws.Range("M21").Value = 1.049972427185191
$ws.Range("N21").Value = 1.038933551196742
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030158046189369
$ws.Range("D22").Value = 1.034785359668948
$ws.Range("E22").Value = 1.03843440920567
$ws.Range("F22").Value = 1.045689378092637
$ws.Range("I22").Value = 1.035710901673691
$ws.Range("J22").Value = 1.037031752090213
$ws.Range("K22").Value = 1.038499178952372
$ws.Range("L22").Value = 1.042134144088949
$ws.Range("M22").Value = 1.049361434333759
$ws.Range("N22").Value = 1.038504456127386
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03051795407659
$ws.Range("D23").Value = 1.034980042152478
$ws.Range("E23").Value = 1.038766209568323
$ws.Range("F23").Value = 1.046082237957161
$ws.Range("I23").Value = 1.035777235499511
$ws.Range("J23").Value = 1.037258916638442
$ws.Range("K23").Value = 1.038624741767799
$ws.Range("L23").Value = 1.042396590467977
$ws.Range("M23").Value = 1.04968527573586
$ws.Range("N23").Value = 1.038731943275329
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031935879591159
$ws.Range("D24").Value = 1.035746913804287
$ws.Range("E24").Value = 1.040074347244037
$ws.Range("F24").Value = 1.047631511570329
$ws.Range("I24").Value = 1.036036005775781
$ws.Range("J24").Value = 1.038152889528635
$ws.Range("K24").Value = 1.039117858658393
$ws.Range("L24").Value = 1.043430251884044
$ws.Range("M24").Value = 1.050961478578581
$ws.Range("N24").Value = 1.039627185709566
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033583352339339
$ws.Range("D25").Value = 1.036637640547429
$ws.Range("E25").Value = 1.041596202846229
$ws.Range("F25").Value = 1.049434750499558
$ws.Range("I25").Value = 1.036331343508695
$ws.Range("J25").Value = 1.039189535959773
$ws.Range("K25").Value = 1.039687539689703
$ws.Range("L25").Value = 1.044630627038571
$ws.Range("M25").Value = 1.040665304297541
